$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row with new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style (format) of an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the team record data for each player row (2 through 42)
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95   # AD
    $ws.Cells.Item($r, 31).Value = 67   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
